$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like numbers to Excel (e.g. "578.99") must be
# forced to remain plain text so they match the original inline-string type
# and exact formatting (trailing zeros, etc).
$numericLookingCells = @("D5", "D6", "D10", "D11", "D12", "D13", "D19", "D20", "D23", "D24", "D25", "D26", "D28", "D31", "D34", "D36", "D37", "D41", "D42", "D46", "D48")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated price values (column D) that need to stay textual
$ws.Range("D5").Value = "578.99"
$ws.Range("D6").Value = "173.21"
$ws.Range("D10").Value = "0.153"
$ws.Range("D11").Value = "0.479"
$ws.Range("D12").Value = "0.0000247"
$ws.Range("D13").Value = "36.80"
$ws.Range("D19").Value = "16.47"
$ws.Range("D20").Value = "491.46"
$ws.Range("D23").Value = "83.90"
$ws.Range("D24").Value = "13.16"
$ws.Range("D25").Value = "2.28"
$ws.Range("D26").Value = "10.49"
$ws.Range("D28").Value = "7.91"
$ws.Range("D31").Value = "28.36"
$ws.Range("D34").Value = "1.00"
$ws.Range("D36").Value = "0.973"
$ws.Range("D37").Value = "46.97"
$ws.Range("D41").Value = "8.45"
$ws.Range("D42").Value = "383.48"
$ws.Range("D46").Value = "135.65"
$ws.Range("D48").Value = "24.87"

# Apply the remaining updated price (column D, already text-safe) and volume
# percentage (column E) values
$ws.Range("D2").Value = "67.055.68"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "3.108.38"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("E9").Value = "  +1.47%  "
$ws.Range("E10").Value = "  -0.89%  "
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("D15").Value = "3.624.20"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Value = "67.042.66"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("D18").Value = "3.110.24"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("E21").Value = "  +4.72%  "
$ws.Range("E22").Value = "  -1.69%  "
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("E25").Value = "  -3.08%  "
$ws.Range("E26").Value = "  +4.77%  "
$ws.Range("E28").Value = "  -0.95%  "
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("E31").Value = "  -1.47%  "
$ws.Range("D33").Value = "0.0₃0944"
$ws.Range("E33").Value = "  -6.25%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("E36").Value = "  -1.65%  "
$ws.Range("E37").Value = "  -1.62%  "
$ws.Range("E38").Value = "  -4.13%  "
$ws.Range("E39").Value = "  -2.60%  "
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("E41").Value = "  -2.41%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Value = "2.803.08"
$ws.Range("E43").Value = "  -1.61%  "
$ws.Range("E44").Value = "  -7.73%  "
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("E51").Value = "  -2.01%  "
